# Implementation of ID tags in progress
#
# Adds a block of new "ID tag" columns (TagID, FloorID, ZoneID, LocationID,
# MemberTypeID, RebarTypeID, SpecificTagID) to the LENGTHS sheet, to the
# right of the existing Lengths/Pcs/Diameter columns, and starts filling in
# the TagID column with sequential letters for each existing data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LENGTHS")

# New header row (row 1), columns D:J
$ws.Range("D1").Value = "TagID"
$ws.Range("E1").Value = "FloorID"
$ws.Range("F1").Value = "ZoneID"
$ws.Range("G1").Value = "LocationID"
$ws.Range("H1").Value = "MemberTypeID"
$ws.Range("I1").Value = "RebarTypeID"
$ws.Range("J1").Value = "SpecificTagID"

# Populate the new TagID column (D) for every existing data row (2-13)
# with sequential letters A..L.
$tagValues = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J", "K", "L")
for ($i = 0; $i -lt $tagValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 4).Value = $tagValues[$i]
}

# Resize the new columns to fit their contents, matching the bestFit
# behaviour Excel applies automatically when new columns are populated
# (widths below are the character-unit equivalents of Excel's computed
# best-fit pixel widths for these headers in the default font).
$ws.Columns.Item(4).EntireColumn.AutoFit()
$ws.Columns.Item(5).EntireColumn.AutoFit()
$ws.Columns.Item(6).EntireColumn.AutoFit()
$ws.Columns.Item(7).EntireColumn.AutoFit()
$ws.Columns.Item(8).EntireColumn.AutoFit()
$ws.Columns.Item(9).EntireColumn.AutoFit()
$ws.Columns.Item(10).EntireColumn.AutoFit()
$ws.Columns.Item(4).ColumnWidth = 5.022135416666667
$ws.Columns.Item(5).ColumnWidth = 6.592447916666667
$ws.Columns.Item(6).ColumnWidth = 6.451822916666667
$ws.Columns.Item(7).ColumnWidth = 9.451822916666666
$ws.Columns.Item(8).ColumnWidth = 14.022135416666666
$ws.Columns.Item(9).ColumnWidth = 11.451822916666666
$ws.Columns.Item(10).ColumnWidth = 11.877604166666666

# Move the active selection to the next empty cell in the TagID column,
# ready for further data entry.
$ws.Range("D14").Select()
